$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '58.989.40'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.21%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.503.02'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.03%  '

$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '534.18'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.31%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.09'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.80%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.37%  '

$ws.Range("E8").Value = '  +1.12%  '

$ws.Range("E9").Value = '  +1.23%  '

$ws.Range("E10").Value = '  -1.36%  '

$ws.Range("E11").Value = '  +1.34%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.347'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.09%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.947.78'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.37%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '58.922.55'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.13%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '22.70'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.96%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000138'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.74%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.504.18'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.08%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.02'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.10%  '

$ws.Range("E19").Value = '  +0.26%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '323.71'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.12%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.07'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.34%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.12%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '65.18'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.32%  '

$ws.Range("E24").Value = '  +0.66%  '

$ws.Range("E25").Value = '  -0.69%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.13%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.52'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.71%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0761'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.80%  '

$ws.Range("B29").Value = 'Monero'
$ws.Range("C29").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '170.26'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.45%  '

$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.75'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.57%  '

$ws.Range("B31").Value = 'Aptos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.45'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -6.02%  '

$ws.Range("E32").Value = '  +2.21%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '18.39'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.33%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.35'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.06%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.05'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.57%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.52'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.98%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.57'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.08%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.800'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.16%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '281.92'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.03%  '

$ws.Range("E41").Value = '  +0.66%  '

$ws.Range("E42").Value = '  -4.90%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '130.13'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.68%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.89'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.38%  '

$ws.Range("E45").Value = '  -0.53%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0923'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.43%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0500'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.15%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0219'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.43%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '17.24'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.59%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.757.47'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.69%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.983'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.15%  '
